$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alloy")

# Select the "alloy" sheet (moves tabSelected / activeTab onto it, off "currency").
$ws.Activate()

# The edit: every percentage-like coefficient in C2:J25 is divided by 100
# (the numbers were entered as whole percents, e.g. 0.25 meaning 25%, and
# needed to be expressed as true fractions, e.g. 0.0025).
$rng = $ws.Range("C2:J25")
foreach ($cell in $rng.Cells) {
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = $v / 100
    }
}

# Re-normalize the touched cells' display: drop the old "#,##0.00"/"#,##0"
# custom formatting + explicit font back to the sheet's implicit default
# style, keeping the centered alignment.
$rng.Style = "Normal"
$rng.HorizontalAlignment = -4108

# Leave the selection on A26, matching where the editor ended up.
$ws.Range("A26").Select()
